$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing label text (A6): microstrain/K -> ue/K ---
$ws.Range("A6").Value = "Coefficient of thermal expansion (ue/K)"

# --- Add new columns F (HSMVals) and G (HSMVars) ---
# Columns F1:F4 and G1:G4 use the same "Text" number format as columns B:E
# (style index 1 in the original workbook, numFmtId 49 / "@"). Set that
# before assigning values so numeric-looking strings stay text instead of
# becoming numbers. Cells are filled row-by-row (F then G on each row) to
# match the order new shared strings get minted in.
$ws.Range("F1:F4").NumberFormat = "@"
$ws.Range("G1:G4").NumberFormat = "@"
$ws.Range("F5:F6").NumberFormat = "@"
# G5 and G6 keep the default (General) cell style -- no NumberFormat override.

$ws.Range("F1").Value = "HSMVals"
$ws.Range("G1").Value = "HSMVars"

$ws.Range("F2").Value = "1.47"
$ws.Range("G2").Value = "Refractive index"

$ws.Range("F3").Value = "0.527212"
$ws.Range("G3").Value = "Grating period (um)"

$ws.Range("F4").Value = "14"
$ws.Range("G4").Value = "Thermo-optic coeff (ue/K)"

$ws.Range("F5").Value = "20"
$ws.Range("G5").Value = "Ref. temperature (deg C)"

$ws.Range("F6").Value = "0.0015"
$ws.Range("G6").Value = "Fiber length (m)"

# --- Selection / active cell now sits on G6 (was G9) ---
[void]$ws.Range("G6").Select()
